$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-03-05 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-03-06 Thursday", 2)

# Update each multiplication problem in the table, addressing each cell by
# its (row, column) coordinates so that identical old/new values occurring
# at different locations never clash with each other.
$tbl = $d.Tables.Item(1)

function Set-CellText($row, $col, $text) {
    $cell = $tbl.Cell($row, $col)
    $range = $cell.Range
    $range.End = $range.End - 1
    $range.Text = $text
}

Set-CellText 1 1 "99×44=4356"
Set-CellText 1 2 "78×35=2730"
Set-CellText 1 3 "34×78=2652"
Set-CellText 1 4 "99×69=6831"
Set-CellText 1 5 "67×55=3685"

Set-CellText 5 1 "42×25=1050"
Set-CellText 5 2 "13×59=767"
Set-CellText 5 3 "57×58=3306"
Set-CellText 5 4 "33×20=660"
Set-CellText 5 5 "87×83=7221"

Set-CellText 10 1 "72×92=6624"
Set-CellText 10 2 "24×75=1800"
Set-CellText 10 3 "18×56=1008"
Set-CellText 10 4 "46×49=2254"
Set-CellText 10 5 "66×11=726"

Set-CellText 15 1 "15×48=720"
Set-CellText 15 2 "60×48=2880"
Set-CellText 15 3 "98×45=4410"
Set-CellText 15 4 "47×86=4042"
Set-CellText 15 5 "57×66=3762"

Set-CellText 20 1 "47×76=3572"
Set-CellText 20 2 "64×11=704"
Set-CellText 20 3 "83×36=2988"
Set-CellText 20 4 "48×47=2256"
Set-CellText 20 5 "67×16=1072"
